# Updates metrics_2_8 worksheet:
#  - Reorders the model names in column A (rows 2-26)
#  - Sets the metric columns B:I (rows 2-26) to the new shared values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "model_2_8_0",
    "model_2_8_22",
    "model_2_8_21",
    "model_2_8_20",
    "model_2_8_19",
    "model_2_8_18",
    "model_2_8_17",
    "model_2_8_16",
    "model_2_8_15",
    "model_2_8_14",
    "model_2_8_13",
    "model_2_8_23",
    "model_2_8_12",
    "model_2_8_10",
    "model_2_8_9",
    "model_2_8_8",
    "model_2_8_7",
    "model_2_8_6",
    "model_2_8_5",
    "model_2_8_4",
    "model_2_8_3",
    "model_2_8_2",
    "model_2_8_1",
    "model_2_8_11",
    "model_2_8_24"
)

$values = @(0.09932080507882668, 0.007716834657778171, -0.6030041075679275, -0.02385691808229473, 0.9967864155769348, 1.138180732727051, 0.2953936159610748, 0.7415750622749329)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
